$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Area / Atotal columns, plus a small Q/A summary block in J:K
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Row 2: per-segment area, running area total, and the J2/K2 summary cells
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Row 3: standalone area formula (not part of the shared group below)
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15 share one formula definition (mirrors the existing D/E shared formulas)
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Reflect the new active selection / scroll position from the saved workbook
$ws.Range("B1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("J2:K2").Select() | Out-Null
